$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.00581
$ws.Range("H2").Value = 3.01743
$ws.Range("I2").Value = 0.003799625168827527
$ws.Range("J2").Value = 0.003799625168827527
$ws.Range("M2").Value = 11.61289466666667
$ws.Range("N2").Value = 34.838684
$ws.Range("O2").Value = 0.09693042549509606
$ws.Range("P2").Value = 0.09693042549509606
$ws.Range("Q2").Value = 11.68036558468
$ws.Range("R2").Value = 105.12329026212
$ws.Range("S2").Value = 0.0003682992843363284
$ws.Range("T2").Value = 0.0003682992843363284

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.00581
$ws.Range("H3").Value = 3.01743
$ws.Range("I3").Value = 0.003799625168827527
$ws.Range("J3").Value = 0.003799625168827527
$ws.Range("O3").Value = 0.2981108740043866
$ws.Range("P3").Value = 0.2981108740043866
$ws.Range("Q3").Value = 35.92312708166001
$ws.Range("R3").Value = 323.30814373494
$ws.Range("S3").Value = 0.001132709579968239
$ws.Range("T3").Value = 0.001132709579968239

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.00581
$ws.Range("H4").Value = 3.01743
$ws.Range("I4").Value = 0.003799625168827527
$ws.Range("J4").Value = 0.003799625168827527
$ws.Range("M4").Value = 27.39934733333333
$ws.Range("N4").Value = 82.198042
$ws.Range("O4").Value = 0.2286966748205465
$ws.Range("P4").Value = 0.2286966748205465
$ws.Range("Q4").Value = 27.55853754134001
$ws.Range("R4").Value = 248.02683787206
$ws.Range("S4").Value = 0.0008689616416753132
$ws.Range("T4").Value = 0.0008689616416753132

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.00581
$ws.Range("H5").Value = 3.01743
$ws.Range("I5").Value = 0.003799625168827527
$ws.Range("J5").Value = 0.003799625168827527
$ws.Range("M5").Value = 45.078635
$ws.Range("N5").Value = 135.235905
$ws.Range("O5").Value = 0.3762620256799708
$ws.Range("P5").Value = 0.3762620256799709
$ws.Range("Q5").Value = 45.34054186935001
$ws.Range("R5").Value = 408.06487682415
$ws.Range("S5").Value = 0.001429654662847646
$ws.Range("T5").Value = 0.001429654662847647

# Row 6
$ws.Range("I6").Value = 0.9594121222074437
$ws.Range("J6").Value = 0.9594121222074438
$ws.Range("M6").Value = 11.61289466666667
$ws.Range("N6").Value = 34.838684
$ws.Range("O6").Value = 0.09693042549509606
$ws.Range("P6").Value = 0.09693042549509606
$ws.Range("Q6").Value = 2949.313112697013
$ws.Range("R6").Value = 26543.81801427312
$ws.Range("S6").Value = 0.09299622523072061
$ws.Range("T6").Value = 0.09299622523072062

# Row 7
$ws.Range("I7").Value = 0.9594121222074437
$ws.Range("J7").Value = 0.9594121222074438
$ws.Range("O7").Value = 0.2981108740043866
$ws.Range("P7").Value = 0.2981108740043866
$ws.Range("S7").Value = 0.2860111862816644
$ws.Range("T7").Value = 0.2860111862816644

# Row 8
$ws.Range("I8").Value = 0.9594121222074437
$ws.Range("J8").Value = 0.9594121222074438
$ws.Range("M8").Value = 27.39934733333333
$ws.Range("N8").Value = 82.198042
$ws.Range("O8").Value = 0.2286966748205465
$ws.Range("P8").Value = 0.2286966748205465
$ws.Range("Q8").Value = 6958.579810552541
$ws.Range("R8").Value = 62627.21829497287
$ws.Range("S8").Value = 0.2194143621313662
$ws.Range("T8").Value = 0.2194143621313662

# Row 9
$ws.Range("I9").Value = 0.9594121222074437
$ws.Range("J9").Value = 0.9594121222074438
$ws.Range("M9").Value = 45.078635
$ws.Range("N9").Value = 135.235905
$ws.Range("O9").Value = 0.3762620256799708
$ws.Range("P9").Value = 0.3762620256799709
$ws.Range("Q9").Value = 11448.56757287237
$ws.Range("R9").Value = 103037.1081558514
$ws.Range("S9").Value = 0.3609903485636924
$ws.Range("T9").Value = 0.3609903485636926

# Row 10
$ws.Range("G10").Value = 9.336668333333334
$ws.Range("H10").Value = 28.010005
$ws.Range("I10").Value = 0.03527091597053946
$ws.Range("J10").Value = 0.03527091597053946
$ws.Range("M10").Value = 11.61289466666667
$ws.Range("N10").Value = 34.838684
$ws.Range("O10").Value = 0.09693042549509606
$ws.Range("P10").Value = 0.09693042549509606
$ws.Range("Q10").Value = 108.4257458926022
$ws.Range("R10").Value = 975.83171303342
$ws.Range("S10").Value = 0.003418824892626169
$ws.Range("T10").Value = 0.003418824892626169

# Row 11
$ws.Range("G11").Value = 9.336668333333334
$ws.Range("H11").Value = 28.010005
$ws.Range("I11").Value = 0.03527091597053946
$ws.Range("J11").Value = 0.03527091597053946
$ws.Range("O11").Value = 0.2981108740043866
$ws.Range("P11").Value = 0.2981108740043866
$ws.Range("Q11").Value = 333.4648920349211
$ws.Range("R11").Value = 3001.18402831429
$ws.Range("S11").Value = 0.0105146435869128
$ws.Range("T11").Value = 0.0105146435869128

# Row 12
$ws.Range("G12").Value = 9.336668333333334
$ws.Range("H12").Value = 28.010005
$ws.Range("I12").Value = 0.03527091597053946
$ws.Range("J12").Value = 0.03527091597053946
$ws.Range("M12").Value = 27.39934733333333
$ws.Range("N12").Value = 82.198042
$ws.Range("O12").Value = 0.2286966748205465
$ws.Range("P12").Value = 0.2286966748205465
$ws.Range("Q12").Value = 255.8186186011345
$ws.Range("R12").Value = 2302.36756741021
$ws.Range("S12").Value = 0.008066341200337284
$ws.Range("T12").Value = 0.008066341200337284

# Row 13
$ws.Range("G13").Value = 9.336668333333334
$ws.Range("H13").Value = 28.010005
$ws.Range("I13").Value = 0.03527091597053946
$ws.Range("J13").Value = 0.03527091597053946
$ws.Range("M13").Value = 45.078635
$ws.Range("N13").Value = 135.235905
$ws.Range("O13").Value = 0.3762620256799708
$ws.Range("P13").Value = 0.3762620256799709
$ws.Range("Q13").Value = 420.8842639143917
$ws.Range("R13").Value = 3787.958375229525
$ws.Range("S13").Value = 0.01327110629066321
$ws.Range("T13").Value = 0.01327110629066321

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.4016586666666667
$ws.Range("H14").Value = 1.204976
$ws.Range("I14").Value = 0.001517336653189343
$ws.Range("J14").Value = 0.001517336653189343
$ws.Range("M14").Value = 11.61289466666667
$ws.Range("N14").Value = 34.838684
$ws.Range("O14").Value = 0.09693042549509606
$ws.Range("P14").Value = 0.09693042549509606
$ws.Range("Q14").Value = 4.664419787953778
$ws.Range("R14").Value = 41.979778091584
$ws.Range("S14").Value = 0.000147076087412948
$ws.Range("T14").Value = 0.000147076087412948

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.4016586666666667
$ws.Range("H15").Value = 1.204976
$ws.Range("I15").Value = 0.001517336653189343
$ws.Range("J15").Value = 0.001517336653189343
$ws.Range("O15").Value = 0.2981108740043866
$ws.Range("P15").Value = 0.2981108740043866
$ws.Range("Q15").Value = 14.34548804060089
$ws.Range("R15").Value = 129.109392365408
$ws.Range("S15").Value = 0.0004523345558411657
$ws.Range("T15").Value = 0.0004523345558411657

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.4016586666666667
$ws.Range("H16").Value = 1.204976
$ws.Range("I16").Value = 0.001517336653189343
$ws.Range("J16").Value = 0.001517336653189343
$ws.Range("M16").Value = 27.39934733333333
$ws.Range("N16").Value = 82.198042
$ws.Range("O16").Value = 0.2286966748205465
$ws.Range("P16").Value = 0.2286966748205465
$ws.Range("Q16").Value = 11.00518531744356
$ws.Range("R16").Value = 99.046667856992
$ws.Range("S16").Value = 0.0003470098471677395
$ws.Range("T16").Value = 0.0003470098471677395

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.4016586666666667
$ws.Range("H17").Value = 1.204976
$ws.Range("I17").Value = 0.001517336653189343
$ws.Range("J17").Value = 0.001517336653189343
$ws.Range("M17").Value = 45.078635
$ws.Range("N17").Value = 135.235905
$ws.Range("O17").Value = 0.3762620256799708
$ws.Range("P17").Value = 0.3762620256799709
$ws.Range("Q17").Value = 18.10622442925333
$ws.Range("R17").Value = 162.95601986328
$ws.Range("S17").Value = 0.0005709161627674894
$ws.Range("T17").Value = 0.0005709161627674895
